# Rotate the species-record values in columns A, B, E, F, G, H, Q, R
# across rows 6, 7 and 8 (row 8 -> row 6, row 6 -> row 7, row 7 -> row 8),
# leaving all other columns (C, D, I, P, S, T, U, V, W, Y, Z, AA, AB, AD,
# AE, AG, AT, AW, AX, AY, ...) untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A", "B", "E", "F", "G", "H", "Q", "R")

# Capture current values of the three rows for the columns we need to move.
# NOTE: the COM interop getter for Range.Value must be invoked with "()"
# (Value() ) to actually retrieve the scalar value instead of a property
# descriptor object; the setter uses plain assignment (Value = ...).
$row6 = @{}
$row7 = @{}
$row8 = @{}
foreach ($col in $cols) {
    $row6[$col] = $ws.Range("${col}6").Value()
    $row7[$col] = $ws.Range("${col}7").Value()
    $row8[$col] = $ws.Range("${col}8").Value()
}

# Apply rotation: new row6 = old row8, new row7 = old row6, new row8 = old row7
foreach ($col in $cols) {
    $ws.Range("${col}6").Value = $row8[$col]
    $ws.Range("${col}7").Value = $row6[$col]
    $ws.Range("${col}8").Value = $row7[$col]
}
